# renouvellement-periode-essai-initiative-employeur.docx
#
# 1) "Lettre Recommandée avec Accusé de Réception"
#       -> lower-case the two mis-capitalised words ("Recommandée" ->
#          "recommandée", "Accusé" -> "accusé", "Réception" -> "réception")
#          and italicise the whole line.
# 2) "l'expression de ma considération" -> replace the straight apostrophe
#    with a typographic (curly) apostrophe: l'expression -> l’expression

$d = $word.ActiveDocument

# --- 1) Lettre Recommandée avec Accusé de Réception -------------------
$d.Content.Find.Execute("Recommand", $true, $false, $false, $false, $false, `
    $true, 1, $false, "recommand", 2) | Out-Null
$d.Content.Find.Execute("avec Accus", $true, $false, $false, $false, $false, `
    $true, 1, $false, "avec accus", 2) | Out-Null
$d.Content.Find.Execute("de R", $true, $false, $false, $false, $false, `
    $true, 1, $false, "de r", 2) | Out-Null

$rng = $d.Content.Duplicate
$rng.Find.Execute("Lettre recommandée avec accusé de réception", $true, `
    $false, $false, $false, $false, $true, 1, $false) | Out-Null
$para = $rng.Paragraphs(1)
$para.Range.Italic = $true
$rng.LanguageID = "fr-FR"

# --- 2) l'expression -> l’expression -----------------------------------
$d.Content.Find.Execute("l'expression de ma considération", $true, $false, `
    $false, $false, $false, $true, 1, $false, "l’expression de ma considération", 2) | Out-Null
